$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodic "Actualizar" refresh: the Fecha (last-updated) timestamps in
# column D cascade down one block and a fresh timestamp is stamped on the
# newest block (rows 2-15). Rows 16-29 and 30-43 pick up the timestamp that
# used to belong to the block above them.
$ws.Range("D2:D15").Value = 44263.51489300781
$ws.Range("D16:D29").Value = 44263.49352696759
$ws.Range("D30:D43").Value = 44263.47216412037
